$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2021" column (P) of data, mirroring the existing 2020 (O) column ---

# Row 3: bottom divider row, same border-only style as the rest of the row (copy from A3).
$ws.Range("A3").Copy()
$ws.Range("P3").PasteSpecial(-4122)   # xlPasteFormats

# Row 4: year header "2021", same style as O4 (bold, bottom border).
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P4").Value = 2021

# Row 5: total victims figure, bold/no-border style (matches style id 18, i.e. O8's style).
$ws.Range("O8").Copy()
$ws.Range("P5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P5").Value = 9038

# Row 6: blank sub-header row (male), regular/no-border style (matches O6/O9's style).
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)   # xlPasteFormats

# Row 7: male total, regular/no-border style (matches O6/O9's style, not O7's own numeric style).
$ws.Range("O6").Copy()
$ws.Range("P7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P7").Value = 8587

# Row 8: male minors, regular/no-border style (matches O6/O9's style).
$ws.Range("O6").Copy()
$ws.Range("P8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P8").Value = 451

# Row 9: blank sub-header row (female), regular/no-border style.
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)   # xlPasteFormats

# Rows 10-24: "no data" placeholder text, right-aligned style (same as col O in each of these rows).
$ws.Range("O10:O24").Copy()
$ws.Range("P10:P24").PasteSpecial(-4122)   # xlPasteFormats
for ($r = 10; $r -le 24; $r++) {
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($r, 15).Value()
}

# Row 25: "no data" placeholder text, right-aligned + bottom border style (same as col O25).
$ws.Range("O25").Copy()
$ws.Range("P25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P25").Value = $ws.Range("O25").Value()

# Clear clipboard/marching-ants after all copy/paste operations.
$excel.CutCopyMode = 0

# Move the active selection to Q4, matching the saved view state.
[void]$ws.Range("Q4").Select()
